# Update column C ("Exposure period") values from 24-hour time ranges
# to 12-hour am/pm time ranges, as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "4:35pm - 5:10pm  31/1/2021"
    3  = "8:10am - 9:30am  31/1/2021"
    4  = "11:15am - 12:00pm  1/2/2021"
    5  = "5:19pm - 6:35pm  1/2/2021"
    6  = "10:00pm - 11:15pm  30/1/2021"
    7  = "4:00pm - 5:00pm  31/1/2021"
    8  = "9:45am - 10:25am  5/2/2021"
    9  = "11:00pm - 11:35pm  29/1/2021"
    10 = "5:00pm - 5:45pm  1/2/2021"
    11 = "2:36pm -3:30pm  30/01/2021"
    12 = "11:30am - 12:15pm  1/2/2021"
    13 = "5:00pm - 6:00pm  31/1/2021"
    14 = "6:50pm - 7:30pm  1/2/2021"
    15 = "6:30pm - 7:30pm  1/2/2021"
    16 = "3:45pm - 4:23pm  5/2/2021"
    17 = "3:57pm - 4:30pm  5/2/2021"
    18 = "3:40pm - 4:15pm  5/2/2021"
    19 = "5:44pm - 6:19pm  7/2/2021"
    20 = "6:17pm - 7:02pm  6/2/2021"
    21 = "3:37pm - 4:10pm  5/2/2021"
    22 = "4:30pm - 6:30pm  6/2/2021"
    23 = "3:53pm - 4:28pm  5/2/2021"
    24 = "5:50pm - 6:30pm  5/2/2021"
    25 = "6:50pm - 7:30pm  6/2/2021"
    26 = "1:00pm - 1:52pm  6/2/2021"
    27 = "11:24pm -12:15am  29/1/2021"
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
